$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B34").Value = "Preprocessed"
$ws.Range("C34").Value = "sex"
$ws.Range("D34").Value = "Factor"
$ws.Range("E34").Value = [char]0x201C + "male" + [char]0x201D + " or " + [char]0x201C + "female" + [char]0x201D

[void]$ws.Range("E34").Select()
